# Commit: "Fruta / hortaliza, semanal"
#
# A new weekly price record is inserted into the data table on row 30,
# pushing the former rows 30..156 down to 31..157 (dimension grows from
# A1:R156 to A1:R157).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 30; existing rows 30-156 shift down to 31-157.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A30").Value = 4
$ws.Range("B30").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C30").Value = "Los Lagos"
$ws.Range("D30").Value = 44565
$ws.Range("E30").Value = 10
$ws.Range("F30").Value = 100112039
$ws.Range("G30").Value = "Ciboulette"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 120
$ws.Range("K30").Value = 3500
$ws.Range("L30").Value = 3500
$ws.Range("M30").Value = 3500
$ws.Range("N30").Value = "$/docena de atados"
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 1167
$ws.Range("Q30").Value = 3
$ws.Range("R30").Value = "Hortaliza"
